$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '37.855.14'
$ws.Range("E2").Value = '  +1.37%  '

$ws.Range("D3").Value = '2.101.83'
$ws.Range("E3").Value = '  +2.12%  '

$ws.Range("E4").Value = '  +0.14%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '235.27'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.08%  '

$ws.Range("E6").Value = '  +1.15%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '58.30'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  +0.71%  '

$ws.Range("E8").Value = '  +0.03%  '

$ws.Range("E9").Value = '  +2.49%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0778'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.25%  '

$ws.Range("D12").Value = '2.417.56'
$ws.Range("E12").Value = '  +2.46%  '

$ws.Range("E13").Value = '  +0.99%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.15'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +1.74%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.784'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +1.43%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.23'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.05%  '

$ws.Range("D17").Value = '2.103.71'
$ws.Range("E17").Value = '  +2.38%  '

$ws.Range("D18").Value = '37.827.73'
$ws.Range("E18").Value = '  +0.75%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.30'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +1.68%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '70.11'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.44%  '

$ws.Range("D21").Value = '0.0₃0822'
$ws.Range("E21").Value = '  +1.33%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '227.51'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.70%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +1.02%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.41'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.08%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '167.85'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.37%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.95'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.59%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.133'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +3.07%  '

$ws.Range("E29").Value = '  -3.27%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '19.51'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +2.19%  '

$ws.Range("E31").Value = '  +0.71%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '4.65'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.98%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '2.58'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.89%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0621'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.36%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.59'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.58%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.47'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +5.47%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.78'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.62%  '

$ws.Range("E38").Value = '  +0.05%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '5.55'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -6.28%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.0967'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.33%  '

$ws.Range("E41").Value = '  -0.50%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '96.82'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.72%  '

$ws.Range("D43").Value = '1.471.61'
$ws.Range("E43").Value = '  +0.79%  '

$ws.Range("E44").Value = '  +0.94%  '

$ws.Range("E45").Value = '  -0.08%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '4.22'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -10.57%  '

$ws.Range("E47").Value = '  +2.15%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '15.50'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.82%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.37'
$ws.Range("D49").Style = "Normal"

$ws.Range("E50").Value = '  +3.07%  '

$ws.Range("D51").Value = '2.304.38'
$ws.Range("E51").Value = '  +2.55%  '

